# Update cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.954.55"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.675.14"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'214.90"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("E6").Value = "  +1.33%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").Value = "'0.0620"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").Value = "'0.0886"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "1.911.33"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").Value = "1.675.98"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "'65.82"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "26.964.15"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").Value = "'237.03"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("E19").Value = "  +3.63%  "

$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").Value = "'4.44"
$ws.Range("E22").Value = "  -0.94%  "

$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("E24").Value = "  -2.41%  "

$ws.Range("D25").Value = "'145.69"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("D27").Value = "'16.00"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "1.483.98"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("D34").Value = "'3.14"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("E35").Value = "  +3.53%  "

$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("E38").Value = "  +1.39%  "

$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("E40").Value = "  -3.16%  "

$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  +2.17%  "

$ws.Range("D44").Value = "'67.30"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "1.817.90"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "'90.52"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").Value = "'7.68"
$ws.Range("E51").Value = "  +0.48%  "
